$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("info")
$ws.Range("E2").Value = " 0 // 0 - expertos; 1 - Igual importancia; 2 - Enfoque Ambiental; 3 - Enfoque Económico; 4 - Enfoque Técnico"

$ws = $wb.Worksheets.Item("alternative_info")
$ws.Range("G4").Value = 640.793947259624
$ws.Range("J4").Value = 223.5449862521488
$ws.Range("L4").Value = 0.0007506139335251635
$ws.Range("M4").Value = 183.2282421272689
$ws.Range("G6").Value = 1067.989912099374
$ws.Range("J6").Value = 278.00170466371
$ws.Range("L6").Value = 0.0006455764403994414
$ws.Range("M6").Value = 175.0574402764045
$ws.Range("G7").Value = 640.793947259624
$ws.Range("J7").Value = 1023.077989380044
$ws.Range("L7").Value = 0.6537430848820482
$ws.Range("M7").Value = 122.2633983069356
$ws.Range("G9").Value = 640.793947259624
$ws.Range("J9").Value = 153.5523892969807
$ws.Range("L9").Value = 0.0005999811623896996
$ws.Range("M9").Value = 171.1756420043254
$ws.Range("G11").Value = 1067.989912099374
$ws.Range("J11").Value = 208.0091077085419
$ws.Range("L11").Value = 0.0004769226130545392
$ws.Range("M11").Value = 161.6726742933793
$ws.Range("G12").Value = 640.793947259624
$ws.Range("J12").Value = 953.0853924248759
$ws.Range("L12").Value = 0.7506442802556406
$ws.Range("M12").Value = 105.8177015515107
$ws.Range("G14").Value = 1067.989912099374
$ws.Range("J14").Value = 1007.542110836437
$ws.Range("L14").Value = 0.7043540549013443
$ws.Range("M14").Value = 107.4352358608177
$ws.Range("G15").Value = 640.793947259624
$ws.Range("J15").Value = 1485.968526752771
$ws.Range("L15").Value = 0.8341918991733863
$ws.Range("M15").Value = 95.49808888516999
$ws.Range("G17").Value = 1067.989912099374
$ws.Range("J17").Value = 138.0165107533738
$ws.Range("L17").Value = 0.000009280258328686273
$ws.Range("M17").Value = 124.5594727710923
$ws.Range("G18").Value = 640.793947259624
$ws.Range("J18").Value = 883.0927954697078
$ws.Range("L18").Value = 0.8813317657511304
$ws.Range("M18").Value = 81.02873272804024
$ws.Range("G20").Value = 1067.989912099374
$ws.Range("J20").Value = 937.549513881269
$ws.Range("L20").Value = 0.8181862245994608
$ws.Range("M20").Value = 85.51782949900341
$ws.Range("G21").Value = 640.793947259624
$ws.Range("J21").Value = 1415.975929797603
$ws.Range("L21").Value = 0.9257343445581757
$ws.Range("M21").Value = 77.82034940911073
$ws.Range("G22").Value = 2135.979824198747
$ws.Range("J22").Value = 274.3666084183291
$ws.Range("G23").Value = 1067.989912099374
$ws.Range("J23").Value = 1470.432648209164
$ws.Range("L23").Value = 0.8827955896432095
$ws.Range("M23").Value = 81.04021835940149

$ws = $wb.Worksheets.Item("alternatives_norm")
$ws.Range("B2").Value = 0.06014777593876288
$ws.Range("D2").Value = 0.006376773534340698
$ws.Range("E2").Value = 0.1051897125100802
$ws.Range("B3").Value = 0.03025183176213168
$ws.Range("D3").Value = 0.006405320791475658
$ws.Range("E3").Value = 0.1054999228969999
$ws.Range("B4").Value = 0.04774139495238833
$ws.Range("D4").Value = 0.008495410555987133
$ws.Range("E4").Value = 0.1165060582853166
$ws.Range("B5").Value = 0.2010041332104988
$ws.Range("D5").Value = 0.000008913973009086989
$ws.Range("E5").Value = 0.1751527609280545
$ws.Range("B6").Value = 0.05937144644710109
$ws.Range("D6").Value = 0.009877642886712468
$ws.Range("E6").Value = 0.1219439757778351
$ws.Range("B7").Value = 0.2184936964007554
$ws.Range("D7").Value = 0.000009754250074386988
$ws.Range("E7").Value = 0.1746000892532424
$ws.Range("B8").Value = 0.3148093435439518
$ws.Range("D8").Value = 0.000007895736867422827
$ws.Range("E8").Value = 0.1988737897390958
$ws.Range("B9").Value = 0.03279342286407274
$ws.Range("D9").Value = 0.01062828957652983
$ws.Range("E9").Value = 0.1247093336811102
$ws.Range("B10").Value = 0.1860561611221832
$ws.Range("D10").Value = 0.000007654621690542115
$ws.Range("E10").Value = 0.2077011267901769
$ws.Range("B11").Value = 0.04442347435878548
$ws.Range("D11").Value = 0.01337066718958758
$ws.Range("E11").Value = 0.1320396310019464
$ws.Range("B12").Value = 0.2035457243124398
$ws.Range("D12").Value = 0.000008495067107110993
$ws.Range("E12").Value = 0.201735625928372
$ws.Range("B13").Value = 0.2998613714556362
$ws.Range("D13").Value = 0.000007139944581198057
$ws.Range("E13").Value = 0.2309931037030518
$ws.Range("B14").Value = 0.2151757758071526
$ws.Range("D14").Value = 0.000009053363844457265
$ws.Range("E14").Value = 0.1986983142518621
$ws.Range("B15").Value = 0.3173509346458929
$ws.Range("D15").Value = 0.000007644252528296597
$ws.Range("E15").Value = 0.2235353660581025
$ws.Range("B16").Value = 0.0003558875855004723
$ws.Range("D16").Value = 0.118088398784087
$ws.Range("E16").Value = 0.279853175888774
$ws.Range("B17").Value = 0.02947550227046987
$ws.Range("D17").Value = 0.687133192685963
$ws.Range("E17").Value = 0.1713815881031083
$ws.Range("B18").Value = 0.1885977522241242
$ws.Range("D18").Value = 0.000007235383747805779
$ws.Range("E18").Value = 0.2634522290808137
$ws.Range("B19").Value = 0.2849133993673206
$ws.Range("D19").Value = 0.000006383972086262789
$ws.Range("E19").Value = 0.2964581091927996
$ws.Range("B20").Value = 0.200227803718837
$ws.Range("D20").Value = 0.000007793792345333628
$ws.Range("E20").Value = 0.2496228024244284
$ws.Range("B21").Value = 0.3024029625575773
$ws.Range("D21").Value = 0.00000688834066903301
$ws.Range("E21").Value = 0.274313857736245
$ws.Range("B22").Value = 0.05859511695543928
$ws.Range("D22").Value = 0.7164914083528875
$ws.Range("E22").Value = 0.1708185985180097
$ws.Range("B23").Value = 0.31403301405229
$ws.Range("D23").Value = 0.000007223386261952141
$ws.Range("E23").Value = 0.2634148906426186
$ws.Range("B24").Value = 0.3676786146310353
$ws.Range("D24").Value = 0.000006376773534340698
$ws.Range("E24").Value = 0.2964888924554955

$ws = $wb.Worksheets.Item("criteria")
$ws.Range("B2").Value = 0.2771882676036963
$ws.Range("B3").Value = 0.07956211810363849
$ws.Range("B4").Value = 0.07596432089007199
$ws.Range("B5").Value = 0.0887594780761196
$ws.Range("B6").Value = 0.02904115627477501
$ws.Range("B7").Value = 0.03630892285481715
$ws.Range("B8").Value = 0.1080909842885875
$ws.Range("B9").Value = 0.05279605036889801
$ws.Range("B10").Value = 0.252288701539396

$ws = $wb.Worksheets.Item("result")
$ws.Range("A2").Value = 13
$ws.Range("B2").Value = 0.7681419272010025
$ws.Range("A3").Value = 21
$ws.Range("B3").Value = 0.7657534507001874
$ws.Range("A4").Value = 6
$ws.Range("B4").Value = 0.7528274157180608
$ws.Range("A5").Value = 19
$ws.Range("B5").Value = 0.7375326120028111
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = 0.7307310833510421
$ws.Range("A7").Value = 22
$ws.Range("B7").Value = 0.7193489749395756
$ws.Range("A8").Value = 17
$ws.Range("B8").Value = 0.7010309319292853
$ws.Range("A9").Value = 12
$ws.Range("B9").Value = 0.6901790720379353
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = 0.6884519701927975
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 0.6645686717893062
$ws.Range("A12").Value = 18
$ws.Range("B12").Value = 0.6616278908808144
$ws.Range("A13").Value = 3
$ws.Range("B13").Value = 0.6527677930409654
$ws.Range("A14").Value = 16
$ws.Range("B14").Value = 0.6366137515129064
$ws.Range("A15").Value = 8
$ws.Range("B15").Value = 0.6300702150024295
$ws.Range("A16").Value = 4
$ws.Range("B16").Value = 0.5253621510483889
$ws.Range("A17").Value = 2
$ws.Range("B17").Value = 0.5205096505271921
$ws.Range("A18").Value = 9
$ws.Range("B18").Value = 0.5191290736267478
$ws.Range("A19").Value = 0
$ws.Range("B19").Value = 0.5189281343297173
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = 0.5063935595023021
$ws.Range("A21").Value = 1
$ws.Range("B21").Value = 0.5009223423275589
$ws.Range("A22").Value = 14
$ws.Range("B22").Value = 0.4565863830734047
$ws.Range("A23").Value = 15
$ws.Range("B23").Value = 0.4358331301153513
$ws.Range("A24").Value = 20
$ws.Range("B24").Value = 0.2903053774954275
